# Apply scraped-schedule refresh for Línea 141 workbook (commit: Horarios actualizados Línea 141 - 1062)
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:55'
$ws.Cells.Item(3, 1).Value = 'Total filas: 41'
# Row 17
$ws.Cells.Item(17, 1).Value = '05:54:55'
$ws.Cells.Item(17, 2).Value = '05:54'
$ws.Cells.Item(17, 3).Value = '10_OLMOS'
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 'LP1912'
# Row 18
$ws.Cells.Item(18, 1).Value = '05:54:55'
$ws.Cells.Item(18, 2).Value = '05:55'
$ws.Cells.Item(18, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 'LP1912'
# Row 19
$ws.Cells.Item(19, 1).Value = '05:20:30'
$ws.Cells.Item(19, 2).Value = '06:04'
$ws.Cells.Item(19, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(19, 4).Value = 44
$ws.Cells.Item(19, 5).Value = 'LP1912'
# Row 20
$ws.Cells.Item(20, 1).Value = '05:54:55'
$ws.Cells.Item(20, 2).Value = '06:11'
$ws.Cells.Item(20, 3).Value = '215A_EL PATO'
$ws.Cells.Item(20, 4).Value = 17
$ws.Cells.Item(20, 5).Value = 'LP1912'
# Row 21
$ws.Cells.Item(21, 1).Value = '05:54:55'
$ws.Cells.Item(21, 2).Value = '06:13'
$ws.Cells.Item(21, 3).Value = '225_HARAS DEL SUR'
$ws.Cells.Item(21, 4).Value = 19
$ws.Cells.Item(21, 5).Value = 'LP1912'
# Row 22
$ws.Cells.Item(22, 1).Value = '05:20:30'
$ws.Cells.Item(22, 2).Value = '06:14'
$ws.Cells.Item(22, 3).Value = '225_HARAS DEL SUR'
$ws.Cells.Item(22, 4).Value = 54
$ws.Cells.Item(22, 5).Value = 'LP1912'
# Row 23
$ws.Cells.Item(23, 1).Value = '05:54:55'
$ws.Cells.Item(23, 2).Value = '06:20'
$ws.Cells.Item(23, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(23, 4).Value = 26
$ws.Cells.Item(23, 5).Value = 'LP1912'
# Row 24
$ws.Cells.Item(24, 1).Value = '05:20:30'
$ws.Cells.Item(24, 2).Value = '06:21'
$ws.Cells.Item(24, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(24, 4).Value = 61
$ws.Cells.Item(24, 5).Value = 'LP1912'
# Row 25
$ws.Cells.Item(25, 1).Value = '05:54:55'
$ws.Cells.Item(25, 2).Value = '06:26'
$ws.Cells.Item(25, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(25, 4).Value = 32
$ws.Cells.Item(25, 5).Value = 'LP1912'
# Row 26
$ws.Cells.Item(26, 1).Value = '05:20:30'
$ws.Cells.Item(26, 2).Value = '06:27'
$ws.Cells.Item(26, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(26, 4).Value = 67
$ws.Cells.Item(26, 5).Value = 'LP1912'
# Row 27
$ws.Cells.Item(27, 1).Value = '05:54:55'
$ws.Cells.Item(27, 2).Value = '06:29'
$ws.Cells.Item(27, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(27, 4).Value = 35
$ws.Cells.Item(27, 5).Value = 'LP1912'
# Row 28
$ws.Cells.Item(28, 1).Value = '05:54:55'
$ws.Cells.Item(28, 2).Value = '06:31'
$ws.Cells.Item(28, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(28, 4).Value = 37
$ws.Cells.Item(28, 5).Value = 'LP1912'
# Row 29
$ws.Cells.Item(29, 1).Value = '05:54:55'
$ws.Cells.Item(29, 2).Value = '06:43'
$ws.Cells.Item(29, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(29, 4).Value = 49
$ws.Cells.Item(29, 5).Value = 'LP1912'
# Row 30
$ws.Cells.Item(30, 1).Value = '05:20:30'
$ws.Cells.Item(30, 2).Value = '06:44'
$ws.Cells.Item(30, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(30, 4).Value = 84
$ws.Cells.Item(30, 5).Value = 'LP1912'
# Row 31
$ws.Cells.Item(31, 1).Value = '05:54:55'
$ws.Cells.Item(31, 2).Value = '06:46'
$ws.Cells.Item(31, 3).Value = '215C_EL PATO'
$ws.Cells.Item(31, 4).Value = 52
$ws.Cells.Item(31, 5).Value = 'LP1912'
# Row 32
$ws.Cells.Item(32, 1).Value = '05:54:55'
$ws.Cells.Item(32, 2).Value = '06:59'
$ws.Cells.Item(32, 3).Value = '14_ABASTO'
$ws.Cells.Item(32, 4).Value = 65
$ws.Cells.Item(32, 5).Value = 'LP1912'
# Row 33
$ws.Cells.Item(33, 1).Value = '05:54:55'
$ws.Cells.Item(33, 2).Value = '07:04'
$ws.Cells.Item(33, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(33, 4).Value = 70
$ws.Cells.Item(33, 5).Value = 'LP1912'
# Row 34
$ws.Cells.Item(34, 1).Value = '05:54:55'
$ws.Cells.Item(34, 2).Value = '07:05'
$ws.Cells.Item(34, 3).Value = '15_ABASTO'
$ws.Cells.Item(34, 4).Value = 71
$ws.Cells.Item(34, 5).Value = 'LP1912'
# Row 35
$ws.Cells.Item(35, 1).Value = '05:54:55'
$ws.Cells.Item(35, 2).Value = '07:06'
$ws.Cells.Item(35, 3).Value = '225_GOMEZ'
$ws.Cells.Item(35, 4).Value = 72
$ws.Cells.Item(35, 5).Value = 'LP1912'
# Row 36
$ws.Cells.Item(36, 1).Value = '05:20:30'
$ws.Cells.Item(36, 2).Value = '07:07'
$ws.Cells.Item(36, 3).Value = '225_GOMEZ'
$ws.Cells.Item(36, 4).Value = 107
$ws.Cells.Item(36, 5).Value = 'LP1912'
# Row 37
$ws.Cells.Item(37, 1).Value = '05:54:55'
$ws.Cells.Item(37, 2).Value = '07:11'
$ws.Cells.Item(37, 3).Value = '215A_EL PATO'
$ws.Cells.Item(37, 4).Value = 77
$ws.Cells.Item(37, 5).Value = 'LP1912'
# Row 38
$ws.Cells.Item(38, 1).Value = '05:54:55'
$ws.Cells.Item(38, 2).Value = '07:15'
$ws.Cells.Item(38, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(38, 4).Value = 81
$ws.Cells.Item(38, 5).Value = 'LP1912'
# Row 39
$ws.Cells.Item(39, 1).Value = '05:54:55'
$ws.Cells.Item(39, 2).Value = '07:20'
$ws.Cells.Item(39, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(39, 4).Value = 86
$ws.Cells.Item(39, 5).Value = 'LP1912'
# Row 40
$ws.Cells.Item(40, 1).Value = '05:54:55'
$ws.Cells.Item(40, 2).Value = '07:22'
$ws.Cells.Item(40, 3).Value = '10_OLMOS'
$ws.Cells.Item(40, 4).Value = 88
$ws.Cells.Item(40, 5).Value = 'LP1912'
# Row 41
$ws.Cells.Item(41, 1).Value = '05:54:55'
$ws.Cells.Item(41, 2).Value = '07:31'
$ws.Cells.Item(41, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(41, 4).Value = 97
$ws.Cells.Item(41, 5).Value = 'LP1912'
# Row 42
$ws.Cells.Item(42, 1).Value = '05:54:55'
$ws.Cells.Item(42, 2).Value = '07:31'
$ws.Cells.Item(42, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(42, 4).Value = 97
$ws.Cells.Item(42, 5).Value = 'LP1912'
# Row 43
$ws.Cells.Item(43, 1).Value = '05:54:55'
$ws.Cells.Item(43, 2).Value = '07:32'
$ws.Cells.Item(43, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(43, 4).Value = 98
$ws.Cells.Item(43, 5).Value = 'LP1912'
# Row 44
$ws.Cells.Item(44, 1).Value = '05:54:55'
$ws.Cells.Item(44, 2).Value = '07:36'
$ws.Cells.Item(44, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(44, 4).Value = 102
$ws.Cells.Item(44, 5).Value = 'LP1912'
# Row 45
$ws.Cells.Item(45, 1).Value = '05:54:55'
$ws.Cells.Item(45, 2).Value = '07:47'
$ws.Cells.Item(45, 3).Value = '14_ABASTO'
$ws.Cells.Item(45, 4).Value = 113
$ws.Cells.Item(45, 5).Value = 'LP1912'
# Row 46
$ws.Cells.Item(46, 1).Value = '05:54:55'
$ws.Cells.Item(46, 2).Value = '07:51'
$ws.Cells.Item(46, 3).Value = '215D_EL PATO'
$ws.Cells.Item(46, 4).Value = 117
$ws.Cells.Item(46, 5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:55'
$ws.Cells.Item(3, 1).Value = 'Total filas: 7'
# Row 9
$ws.Cells.Item(9, 1).Value = '05:54:55'
$ws.Cells.Item(9, 2).Value = '06:11'
$ws.Cells.Item(9, 3).Value = '215A_EL PATO'
$ws.Cells.Item(9, 4).Value = 17
$ws.Cells.Item(9, 5).Value = 'LP1912'
# Row 10
$ws.Cells.Item(10, 1).Value = '05:54:55'
$ws.Cells.Item(10, 2).Value = '06:46'
$ws.Cells.Item(10, 3).Value = '215C_EL PATO'
$ws.Cells.Item(10, 4).Value = 52
$ws.Cells.Item(10, 5).Value = 'LP1912'
# Row 11
$ws.Cells.Item(11, 1).Value = '05:54:55'
$ws.Cells.Item(11, 2).Value = '07:11'
$ws.Cells.Item(11, 3).Value = '215A_EL PATO'
$ws.Cells.Item(11, 4).Value = 77
$ws.Cells.Item(11, 5).Value = 'LP1912'
# Row 12
$ws.Cells.Item(12, 1).Value = '05:54:55'
$ws.Cells.Item(12, 2).Value = '07:51'
$ws.Cells.Item(12, 3).Value = '215D_EL PATO'
$ws.Cells.Item(12, 4).Value = 117
$ws.Cells.Item(12, 5).Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:55'
$ws.Cells.Item(3, 1).Value = 'Total filas: 8'
# Row 7
$ws.Cells.Item(7, 1).Value = '05:54:55'
$ws.Cells.Item(7, 2).Value = '06:09'
$ws.Cells.Item(7, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(7, 4).Value = 15
$ws.Cells.Item(7, 5).Value = 'L6173'
# Row 9
$ws.Cells.Item(9, 1).Value = '05:54:55'
$ws.Cells.Item(9, 2).Value = '06:32'
$ws.Cells.Item(9, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(9, 4).Value = 38
$ws.Cells.Item(9, 5).Value = 'L6203'
# Row 10
$ws.Cells.Item(10, 1).Value = '05:20:30'
$ws.Cells.Item(10, 2).Value = '06:33'
$ws.Cells.Item(10, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(10, 4).Value = 73
$ws.Cells.Item(10, 5).Value = 'L6203'
# Row 11
$ws.Cells.Item(11, 1).Value = '05:54:55'
$ws.Cells.Item(11, 2).Value = '06:59'
$ws.Cells.Item(11, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(11, 4).Value = 65
$ws.Cells.Item(11, 5).Value = 'L6173'
# Row 12
$ws.Cells.Item(12, 1).Value = '05:20:30'
$ws.Cells.Item(12, 2).Value = '07:00'
$ws.Cells.Item(12, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(12, 4).Value = 100
$ws.Cells.Item(12, 5).Value = 'L6173'
# Row 13
$ws.Cells.Item(13, 1).Value = '05:54:55'
$ws.Cells.Item(13, 2).Value = '07:34'
$ws.Cells.Item(13, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(13, 4).Value = 100
$ws.Cells.Item(13, 5).Value = 'L6173'
